# Generate Report for Archive
# - Flip the localization status shown on the Overview / zh-cn / de-de
#   sheets from "Ready for handoff" to "In Translation".
# - Re-fit the now-shorter status columns (Overview!E:F and the "Status"
#   column on each language sheet) to their new, narrower width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update every cell currently showing the old "Ready for handoff" status.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Narrow the affected columns to match the new content width.
$newColumnWidth = 12.576851254417766

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
